$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab (Property1 -> DataNode), per the commit message:
# "unify the conception of DataNode, DataTable, Entity."
$ws.Name = "DataNode"

# Column A was nudged slightly narrower in the saved view state.
$ws.Columns.Item(1).ColumnWidth = 23.38

# Move/restore the active selection to C24, matching the saved view state.
$ws.Range("C24").Select() | Out-Null
